$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (D and E) keep their original text formatting
# instead of being auto-converted to numbers by Excel type inference.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.579.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.580.62"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.23%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.03"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.503"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.17%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.10"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.252"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0589"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.802.17"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.543.21"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.85"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.527"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.552.81"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.76"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.69"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.97%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.32"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.83%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0691"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.96%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.44"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -6.24%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.26"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.68"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.03"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.37%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.74%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.365.58"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.94"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.51"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.964"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.16%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0164"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.35%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.37%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.975"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.98%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.41"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.76"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.72%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.24"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.712.87"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.29"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0998"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0968"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.79%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.79%  "
